$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new value map for this update pass
$updates = [ordered]@{
    'D2' = '43.724.01'
    'E2' = '  -0.58%  '
    'D3' = '2.229.89'
    'E3' = '  -1.85%  '
    'E4' = '  -0.02%  '
    'D5' = '312.86'
    'E5' = '  -1.63%  '
    'D6' = '98.11'
    'E6' = '  -4.29%  '
    'D7' = '0.567'
    'E7' = '  -3.22%  '
    'E8' = '  +0.03%  '
    'D9' = '0.532'
    'E9' = '  -6.58%  '
    'D10' = '35.79'
    'E10' = '  -6.82%  '
    'D11' = '0.0818'
    'E11' = '  -2.05%  '
    'D12' = '7.35'
    'E12' = '  -6.32%  '
    'E13' = '  -2.97%  '
    'D14' = '2.564.94'
    'E14' = '  -2.42%  '
    'B15' = 'Polygon'
    'C15' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    'D15' = '0.838'
    'E15' = '  -4.00%  '
    'B16' = 'WrappedEther'
    'C16' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    'D16' = '2.226.22'
    'E16' = '  -2.49%  '
    'D17' = '14.00'
    'E17' = '  -3.37%  '
    'D18' = '43.587.22'
    'E18' = '  -0.86%  '
    'D19' = '13.01'
    'E19' = '  -8.96%  '
    'D20' = '0.0₃0962'
    'E20' = '  -3.35%  '
    'D21' = '6.29'
    'E21' = '  -5.30%  '
    'D22' = '65.06'
    'E22' = '  -1.56%  '
    'D23' = '235.15'
    'E23' = '  -1.04%  '
    'D24' = '2.96'
    'E24' = '  -7.66%  '
    'D25' = '2.02'
    'E25' = '  -8.11%  '
    'E26' = '  +0.28%  '
    'D27' = '9.98'
    'E27' = '  -2.88%  '
    'D28' = '2.19'
    'E28' = '  -1.55%  '
    'D29' = '36.57'
    'E29' = '  -6.85%  '
    'D30' = '5.97'
    'E30' = '  -7.99%  '
    'D31' = '157.25'
    'E31' = '  -2.96%  '
    'D32' = '19.83'
    'E32' = '  -3.04%  '
    'D33' = '0.0826'
    'E33' = '  -5.76%  '
    'D34' = '2.66'
    'E34' = '  -2.38%  '
    'D35' = '3.13'
    'E35' = '  -4.87%  '
    'D36' = '0.109'
    'E36' = '  +0.13%  '
    'D37' = '1.88'
    'E37' = '  -7.27%  '
    'D38' = '0.117'
    'E38' = '  -3.27%  '
    'D39' = '15.42'
    'E39' = '  -0.44%  '
    'B40' = 'RenderToken'
    'C40' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D40' = '4.02'
    'E40' = '  -11.14%  '
    'B41' = 'NEARProtocol'
    'C41' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D41' = '3.53'
    'E41' = '  -8.25%  '
    'D42' = '0.0306'
    'E42' = '  -6.16%  '
    'E43' = '  +0.00%  '
    'D44' = '1.702.15'
    'E44' = '  -4.10%  '
    'D45' = '83.28'
    'E45' = '  -1.81%  '
    'D46' = '0.193'
    'E46' = '  -6.64%  '
    'D47' = '5.11'
    'E47' = '  -5.58%  '
    'D48' = '101.37'
    'E48' = '  -2.72%  '
    'D49' = '1.63'
    'E49' = '  +1.43%  '
    'D50' = '70.90'
    'E50' = '  -4.78%  '
    'D51' = '55.92'
    'E51' = '  -5.78%  '
}

# Force text storage (column D holds numeric-looking strings such as
# "43.724.01" / "2.226.22" that Excel would otherwise reinterpret as
# numbers/dates) then restore the default "Normal" style so no visible
# formatting change is introduced.
foreach ($key in $updates.Keys) {
    $cell = $ws.Range($key)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$key]
    $cell.Style = "Normal"
}
